# Weekly fruit/vegetable price update: insert a new current-week row at
# row 7 (pushing the existing rows 7-17 down to 8-18) and populate it
# with the latest "Espinaca" price observation for
# "Terminal Hortofrutícola Agro Chillán".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 7:17 down to 8:18, duplicating row 7's
# formatting (incl. the date number format on column D) into the new row.
$ws.Rows.Item(7).Insert()

# Fill the newly inserted row 7 with this week's observation.
$ws.Cells.Item(7, 1).Value  = 7
$ws.Cells.Item(7, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value  = "Ñuble"
$ws.Cells.Item(7, 4).Value  = 44819
$ws.Cells.Item(7, 5).Value  = 16
$ws.Cells.Item(7, 6).Value  = 100112012
$ws.Cells.Item(7, 7).Value  = "Espinaca"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 7000
$ws.Cells.Item(7, 12).Value = 8000
$ws.Cells.Item(7, 13).Value = 7500
$ws.Cells.Item(7, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(7, 16).Value = 750
$ws.Cells.Item(7, 17).Value = 10
$ws.Cells.Item(7, 18).Value = "Hortaliza"
